$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row appended to the daily log: 2025/10/10, 金, 20, 201
# Column A holds date-like text ("2025/10/10") that must stay literal text
# (matching every other row in the sheet) rather than Excel's automatic
# date-serial conversion, so the cell is pre-formatted as Text before the
# value is typed in, then the formatting is reset back to Normal/General
# (only the *value* differs from the rest of the sheet's default style).
$dateCell = $ws.Range("A90")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/10/10"
$dateCell.Style = "Normal"

$ws.Range("B90").Value = "金"
$ws.Range("C90").Value = 20
$ws.Range("D90").Value = 201
